$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LocationZ values for rows 2 through 62 (column B)
$newValues = @(2923.576028960172, 3321.038713266705, 3886.777020343889, 4314.744431086633, 4642.171555175671, 4868.624648148253, 5125.944010635396, 5342.756567636689, 5626.781153722075, 5727.118723810588, 5791.369956989529, 6019.269040680137, 6162.996299292629, 6220.965503800751, 6381.109695094492, 6441.654000796171, 6468.345827178451, 6498.077576104009, 6597.705056959815, 6615.022245438043, 6627.878116110614, 6602.238769319099, 6628.409705467133, 6646.820530222499, 6618.202989504023, 6600.240194404897, 6508.554519089291, 6494.293126123084, 6434.584857169063, 6331.740538853276, 6299.064790382066, 6216.693799803349, 6144.14444510644, 6057.392124222489, 5973.542705715726, 5840.625810158159, 5693.227347891461, 5591.618838359409, 5514.104039840689, 5308.078403325559, 5189.985974329237, 5047.705473335253, 4948.090802797249, 4766.334184526131, 4583.11242521156, 4448.042991517098, 4271.029509850966, 4091.851195416982, 3970.624824456323, 3783.0142072356, 3590.269155933204, 3442.198661563002, 3302.998262307523, 3058.565881269529, 2866.532324447465, 2688.506866639567, 2574.799292610603, 2460.108900904722, 2382.753069565312, 2264.021341908573, 2251.988430747316)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
